$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) "Ci-dessous ... rattachées. " paragraph: the sentence was split
#    across two runs ("...rattac" / "hées. "). Re-typing it as one
#    contiguous string via Find/Replace consolidates it into a single
#    run, matching the target OOXML.
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "Ci-dessous on peut retrouver la liste complète des fonctionnalités répartie selon leur fonctionnalité principale à laquelle elles sont rattachées. ",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Ci-dessous on peut retrouver la liste complète des fonctionnalités répartie selon leur fonctionnalité principale à laquelle elles sont rattachées. ",
    2) | Out-Null

# ------------------------------------------------------------------
# 2) "Gestion client" paragraph: was split across several runs with
#    the _GoBack bookmark sitting in the middle of the sentence.
#    Remove that stray bookmark first (Find/Replace across a bookmark
#    would just delete it anyway), then retype the sentence as one
#    contiguous run.
# ------------------------------------------------------------------
try {
    $bm = $d.Bookmarks.Item("_GoBack")
    $bm.Delete() | Out-Null
} catch {
}

$d.Content.Find.Execute(
    "En ce qu’il s’agit de la gestion client on pourra retrouver la liste des clients ainsi que leur info pour chacun. La possibilité d’ajouter et d’éditer un client sera présente également et l’affichage des erreurs dans les champs des formulaires liés aux clients.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "En ce qu’il s’agit de la gestion client on pourra retrouver la liste des clients ainsi que leur info pour chacun. La possibilité d’ajouter et d’éditer un client sera présente également et l’affichage des erreurs dans les champs des formulaires liés aux clients.",
    2) | Out-Null

# ------------------------------------------------------------------
# 3) "Dans la gestion de projets ..." paragraph: same text, just
#    split across runs -> retype to consolidate into a single run.
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "Dans la gestion de projets nous pourront retrouver la liste des projets ainsi que les infos de chacun. On pourra également créer et éditer un projet. Ces dernières fonctionnalités seront également présentes pour les plans d’un projet. On pourra récupérer la liste des clients, consulter le devis, copier le plan charger les données paramètres. Et enfin pour tous les champs présents dans les fonctionnalités précédentes on affichera les erreurs éventuelles.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Dans la gestion de projets nous pourront retrouver la liste des projets ainsi que les infos de chacun. On pourra également créer et éditer un projet. Ces dernières fonctionnalités seront également présentes pour les plans d’un projet. On pourra récupérer la liste des clients, consulter le devis, copier le plan charger les données paramètres. Et enfin pour tous les champs présents dans les fonctionnalités précédentes on affichera les erreurs éventuelles.",
    2) | Out-Null

# ------------------------------------------------------------------
# 4) "Pour la fonctionnalité gestion de devis ..." paragraph: add the
#    new sentence about checking stock levels, and consolidate the
#    runs at the same time.
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "Pour la fonctionnalité gestion de devis, on pourra charger les composants, afficher les infos du plan, charger les données client et les données commerciales, on calculera les prix et la possibilité d’exporter le devis sera présente. On pourra également appliquer une remise, afficher un dossier technique ainsi que la modification de l’état d’un devis. Pour tous les champs présents dans les fonctionnalités précédentes on affichera les erreurs éventuelles. ",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Pour la fonctionnalité gestion de devis, on pourra charger les composants, afficher les infos du plan, charger les données client et les données commerciales, on vérifiera que les stocks sont suffisamment remplis pour le devis, on calculera les prix et la possibilité d’exporter le devis sera présente. On pourra également appliquer une remise, afficher un dossier technique ainsi que la modification de l’état d’un devis. Pour tous les champs présents dans les fonctionnalités précédentes on affichera les erreurs éventuelles. ",
    2) | Out-Null

# ------------------------------------------------------------------
# 5) "En ce qui s’agit de la fonctionnalité de la modélisation ..."
#    paragraph: drop the obsolete slot/module sentence, replace it
#    with the shorter "choisir ou sélectionner un plan" sentence, and
#    consolidate the runs. Note the source has a non-breaking space
#    between "plan" and "et vérifier".
# ------------------------------------------------------------------
$oldModel = "En ce qui s’agit de la fonctionnalité de la modélisation on pourra éditer les paramètres modifiables, sauvegarder le plan" + [char]0x00A0 + "et vérifier ses contraintes. On aura la possibilité de quitté et tracer, slot, le retirer ou en sélectionner un. On pourra également choisir ou sélectionner un module, sélectionner le slot d’un module et choisir un module pour un slot de module. Enfin on pourra mettre à jour l’affichage et charger les modules compatibles et les composants d’affichages. "
$newModel = "En ce qui s’agit de la fonctionnalité de la modélisation on pourra éditer les paramètres modifiables, sauvegarder le plan et vérifier ses contraintes. On pourra également choisir ou sélectionner un plan. Enfin on pourra mettre à jour l’affichage et charger les modules compatibles et les composants d’affichages. "

$d.Content.Find.Execute(
    $oldModel,
    $false, $false, $false, $false, $false, $true, 1, $false,
    $newModel,
    2) | Out-Null

# ------------------------------------------------------------------
# 6) Re-insert the _GoBack bookmark at the very end of the
#    modélisation paragraph (right after the text, before the
#    paragraph mark), matching the target layout.
# ------------------------------------------------------------------
$lastPara = $d.Paragraphs.Last
$bmRange = $lastPara.Range.Duplicate
$bmRange.MoveEnd(1, -1) | Out-Null
$bmRange.Collapse(0) | Out-Null
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
